$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35
$ws.Range("B35").Value = 6847848
$ws.Range("C35").Value = "France Ligue 1"
$ws.Range("D35").Value = 45165.41666666666
$ws.Range("E35").Value = "Clermont Foot"
$ws.Range("F35").Value = "Metz"
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 1
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = "A"
$ws.Range("L35").Value = 1.85
$ws.Range("M35").Value = 3.6
$ws.Range("N35").Value = 4
$ws.Range("O35").Value = 1.909
$ws.Range("P35").Value = 3.6
$ws.Range("Q35").Value = 4
$ws.Range("R35").Value = -0.5
$ws.Range("S35").Value = 1.9
$ws.Range("T35").Value = 1.95
$ws.Range("U35").Value = 2.5
$ws.Range("V35").Value = 1.95
$ws.Range("W35").Value = 1.9
$ws.Range("X35").Value = -1
$ws.Range("Y35").Value = -1
$ws.Range("Z35").Value = 3
$ws.Range("AA35").Value = -1
$ws.Range("AB35").Value = 0.95
$ws.Range("AC35").Value = -1
$ws.Range("AD35").Value = 0.8999999999999999

# Row 36
$ws.Range("B36").Value = 6847850
$ws.Range("C36").Value = "France Ligue 1"
$ws.Range("D36").Value = 45165.41666666666
$ws.Range("E36").Value = "Montpellier"
$ws.Range("F36").Value = "Reims"
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 3
$ws.Range("I36").Value = 1
$ws.Range("J36").Value = 2
$ws.Range("K36").Value = "A"
$ws.Range("L36").Value = 2.4
$ws.Range("M36").Value = 3.4
$ws.Range("N36").Value = 2.8
$ws.Range("O36").Value = 2.2
$ws.Range("P36").Value = 3.5
$ws.Range("Q36").Value = 3.2
$ws.Range("R36").Value = -0.25
$ws.Range("S36").Value = 1.95
$ws.Range("T36").Value = 1.9
$ws.Range("U36").Value = 2.75
$ws.Range("V36").Value = 1.925
$ws.Range("W36").Value = 1.925
$ws.Range("X36").Value = -1
$ws.Range("Y36").Value = -1
$ws.Range("Z36").Value = 2.2
$ws.Range("AA36").Value = -1
$ws.Range("AB36").Value = 0.8999999999999999
$ws.Range("AC36").Value = 0.925
$ws.Range("AD36").Value = -1

# Row 121
$ws.Range("B121").Value = 6847939
$ws.Range("C121").Value = "France Ligue 1"
$ws.Range("D121").Value = 45256.45833333334
$ws.Range("E121").Value = "Nantes"
$ws.Range("F121").Value = "Le Havre"
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = "D"
$ws.Range("L121").Value = 2.15
$ws.Range("M121").Value = 3.3
$ws.Range("N121").Value = 3.4
$ws.Range("O121").Value = 1.95
$ws.Range("P121").Value = 3.25
$ws.Range("Q121").Value = 4.333
$ws.Range("R121").Value = -0.5
$ws.Range("S121").Value = 1.95
$ws.Range("T121").Value = 1.9
$ws.Range("U121").Value = 2
$ws.Range("V121").Value = 1.825
$ws.Range("W121").Value = 2.025
$ws.Range("X121").Value = -1
$ws.Range("Y121").Value = 2.25
$ws.Range("Z121").Value = -1
$ws.Range("AA121").Value = -1
$ws.Range("AB121").Value = 0.8999999999999999
$ws.Range("AC121").Value = -1
$ws.Range("AD121").Value = 1.025

# Row 122
$ws.Range("B122").Value = 6847940
$ws.Range("C122").Value = "France Ligue 1"
$ws.Range("D122").Value = 45256.45833333334
$ws.Range("E122").Value = "Montpellier"
$ws.Range("F122").Value = "Brest"
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 3
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1
$ws.Range("K122").Value = "A"
$ws.Range("L122").Value = 2.15
$ws.Range("M122").Value = 3.4
$ws.Range("N122").Value = 3.25
$ws.Range("O122").Value = 2.25
$ws.Range("P122").Value = 3.3
$ws.Range("Q122").Value = 3.25
$ws.Range("R122").Value = -0.25
$ws.Range("S122").Value = 1.95
$ws.Range("T122").Value = 1.9
$ws.Range("U122").Value = 2.5
$ws.Range("V122").Value = 2.05
$ws.Range("W122").Value = 1.8
$ws.Range("X122").Value = -1
$ws.Range("Y122").Value = -1
$ws.Range("Z122").Value = 2.25
$ws.Range("AA122").Value = -1
$ws.Range("AB122").Value = 0.8999999999999999
$ws.Range("AC122").Value = 1.05
$ws.Range("AD122").Value = -1

# Row 123
$ws.Range("B123").Value = 6847935
$ws.Range("C123").Value = "France Ligue 1"
$ws.Range("D123").Value = 45256.45833333334
$ws.Range("E123").Value = "Lorient"
$ws.Range("F123").Value = "Metz"
$ws.Range("G123").Value = 2
$ws.Range("H123").Value = 3
$ws.Range("I123").Value = 2
$ws.Range("J123").Value = 1
$ws.Range("K123").Value = "A"
$ws.Range("L123").Value = 1.95
$ws.Range("M123").Value = 3.5
$ws.Range("N123").Value = 3.8
$ws.Range("O123").Value = 2.05
$ws.Range("P123").Value = 3.4
$ws.Range("Q123").Value = 3.8
$ws.Range("R123").Value = -0.5
$ws.Range("S123").Value = 2.05
$ws.Range("T123").Value = 1.8
$ws.Range("U123").Value = 2.25
$ws.Range("V123").Value = 1.975
$ws.Range("W123").Value = 1.875
$ws.Range("X123").Value = -1
$ws.Range("Y123").Value = -1
$ws.Range("Z123").Value = 2.8
$ws.Range("AA123").Value = -1
$ws.Range("AB123").Value = 0.8
$ws.Range("AC123").Value = 0.9750000000000001
$ws.Range("AD123").Value = -1

# Row 196
$ws.Range("B196").Value = 6848007
$ws.Range("C196").Value = "France Ligue 1"
$ws.Range("D196").Value = 45333.45833333334
$ws.Range("E196").Value = "Lorient"
$ws.Range("F196").Value = "Reims"
$ws.Range("G196").Value = 2
$ws.Range("H196").Value = 0
$ws.Range("I196").Value = 0
$ws.Range("J196").Value = 0
$ws.Range("K196").Value = "H"
$ws.Range("L196").Value = 3.3
$ws.Range("M196").Value = 3.5
$ws.Range("N196").Value = 2.1
$ws.Range("O196").Value = 4.2
$ws.Range("P196").Value = 3.75
$ws.Range("Q196").Value = 1.833
$ws.Range("R196").Value = 0.5
$ws.Range("S196").Value = 1.975
$ws.Range("T196").Value = 1.875
$ws.Range("U196").Value = 2.5
$ws.Range("V196").Value = 1.875
$ws.Range("W196").Value = 1.975
$ws.Range("X196").Value = 3.2
$ws.Range("Y196").Value = -1
$ws.Range("Z196").Value = -1
$ws.Range("AA196").Value = 0.9750000000000001
$ws.Range("AB196").Value = -1
$ws.Range("AC196").Value = -1
$ws.Range("AD196").Value = 0.9750000000000001

# Row 197
$ws.Range("B197").Value = 6848005
$ws.Range("C197").Value = "France Ligue 1"
$ws.Range("D197").Value = 45333.45833333334
$ws.Range("E197").Value = "Toulouse"
$ws.Range("F197").Value = "Nantes"
$ws.Range("G197").Value = 1
$ws.Range("H197").Value = 2
$ws.Range("I197").Value = 0
$ws.Range("J197").Value = 1
$ws.Range("K197").Value = "A"
$ws.Range("L197").Value = 2.25
$ws.Range("M197").Value = 3.25
$ws.Range("N197").Value = 3.2
$ws.Range("O197").Value = 2.05
$ws.Range("P197").Value = 3.4
$ws.Range("Q197").Value = 3.6
$ws.Range("R197").Value = -0.5
$ws.Range("S197").Value = 2.05
$ws.Range("T197").Value = 1.8
$ws.Range("U197").Value = 2.5
$ws.Range("V197").Value = 2.025
$ws.Range("W197").Value = 1.825
$ws.Range("X197").Value = -1
$ws.Range("Y197").Value = -1
$ws.Range("Z197").Value = 2.6
$ws.Range("AA197").Value = -1
$ws.Range("AB197").Value = 0.8
$ws.Range("AC197").Value = 1.025
$ws.Range("AD197").Value = -1

# Row 259
$ws.Range("B259").Value = 6848069
$ws.Range("C259").Value = "France Ligue 1"
$ws.Range("D259").Value = 45389.41666666666
$ws.Range("E259").Value = "Reims"
$ws.Range("F259").Value = "Nice"
$ws.Range("G259").Value = 0
$ws.Range("H259").Value = 0
$ws.Range("I259").Value = 0
$ws.Range("J259").Value = 0
$ws.Range("K259").Value = "D"
$ws.Range("L259").Value = 2.45
$ws.Range("M259").Value = 3.1
$ws.Range("N259").Value = 2.9
$ws.Range("O259").Value = 2.75
$ws.Range("P259").Value = 3.1
$ws.Range("Q259").Value = 2.7
$ws.Range("R259").Value = 0
$ws.Range("S259").Value = 1.9
$ws.Range("T259").Value = 1.95
$ws.Range("U259").Value = 2.25
$ws.Range("V259").Value = 1.875
$ws.Range("W259").Value = 1.975
$ws.Range("X259").Value = -1
$ws.Range("Y259").Value = 2.1
$ws.Range("Z259").Value = -1
$ws.Range("AA259").Value = 0
$ws.Range("AB259").Value = 0
$ws.Range("AC259").Value = -1
$ws.Range("AD259").Value = 0.9750000000000001

# Row 260
$ws.Range("B260").Value = 6977742
$ws.Range("C260").Value = "France Ligue 1"
$ws.Range("D260").Value = 45389.41666666666
$ws.Range("E260").Value = "Montpellier"
$ws.Range("F260").Value = "Lorient"
$ws.Range("G260").Value = 2
$ws.Range("H260").Value = 0
$ws.Range("I260").Value = 0
$ws.Range("J260").Value = 0
$ws.Range("K260").Value = "H"
$ws.Range("L260").Value = 1.666
$ws.Range("M260").Value = 3.75
$ws.Range("N260").Value = 5
$ws.Range("O260").Value = 1.8
$ws.Range("P260").Value = 4
$ws.Range("Q260").Value = 4
$ws.Range("R260").Value = -0.5
$ws.Range("S260").Value = 1.825
$ws.Range("T260").Value = 2.025
$ws.Range("U260").Value = 3
$ws.Range("V260").Value = 1.95
$ws.Range("W260").Value = 1.9
$ws.Range("X260").Value = 0.8
$ws.Range("Y260").Value = -1
$ws.Range("Z260").Value = -1
$ws.Range("AA260").Value = 0.825
$ws.Range("AB260").Value = -1
$ws.Range("AC260").Value = -1
$ws.Range("AD260").Value = 0.8999999999999999

# Row 309
$ws.Range("B309").Value = 7162745
$ws.Range("C309").Value = "France Ligue 1"
$ws.Range("D309").Value = 45431.66666666666
$ws.Range("E309").Value = "Lens"
$ws.Range("F309").Value = "Montpellier"
$ws.Range("G309").Value = 2
$ws.Range("H309").Value = 2
$ws.Range("I309").Value = 2
$ws.Range("J309").Value = 0
$ws.Range("K309").Value = "D"
$ws.Range("L309").Value = 1.5
$ws.Range("M309").Value = 4.2
$ws.Range("N309").Value = 6.5
$ws.Range("O309").Value = 1.333
$ws.Range("P309").Value = 5.5
$ws.Range("Q309").Value = 8
$ws.Range("R309").Value = -1.75
$ws.Range("S309").Value = 2
$ws.Range("T309").Value = 1.85
$ws.Range("U309").Value = 3.5
$ws.Range("V309").Value = 1.925
$ws.Range("W309").Value = 1.925
$ws.Range("X309").Value = -1
$ws.Range("Y309").Value = 4.5
$ws.Range("Z309").Value = -1
$ws.Range("AA309").Value = -1
$ws.Range("AB309").Value = 0.8500000000000001
$ws.Range("AC309").Value = 0.925
$ws.Range("AD309").Value = -1

# Row 310
$ws.Range("B310").Value = 7162746
$ws.Range("C310").Value = "France Ligue 1"
$ws.Range("D310").Value = 45431.66666666666
$ws.Range("E310").Value = "Metz"
$ws.Range("F310").Value = "PSG"
$ws.Range("G310").Value = 0
$ws.Range("H310").Value = 2
$ws.Range("I310").Value = 0
$ws.Range("J310").Value = 2
$ws.Range("K310").Value = "A"
$ws.Range("L310").Value = 6
$ws.Range("M310").Value = 4
$ws.Range("N310").Value = 1.533
$ws.Range("O310").Value = 4.5
$ws.Range("P310").Value = 3.6
$ws.Range("Q310").Value = 1.8
$ws.Range("R310").Value = 0.75
$ws.Range("S310").Value = 1.85
$ws.Range("T310").Value = 2
$ws.Range("U310").Value = 3
$ws.Range("V310").Value = 1.825
$ws.Range("W310").Value = 2.025
$ws.Range("X310").Value = -1
$ws.Range("Y310").Value = -1
$ws.Range("Z310").Value = 0.8
$ws.Range("AA310").Value = -1
$ws.Range("AB310").Value = 1
$ws.Range("AC310").Value = -1
$ws.Range("AD310").Value = 1.025
